$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures from the latest data pull.
# D-column price cells are force-typed as text (leading apostrophe + ClearFormats)
# so Excel does not reinterpret numeric-looking strings (e.g. "560.32") as numbers
# and mangle formatting (trailing zeros, float precision, etc.), matching the
# original inline-string representation of these cells.

$ws.Range('D2').Value = "'68.827.04"
$ws.Range('D2').ClearFormats()
$ws.Range('D3').Value = "'2.465.78"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'560.32"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').Value = "'163.92"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.63%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.512"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('D9').Value = "'2.464.71"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('D12').Value = "'0.330"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = "'4.84"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.65%  '
$ws.Range('D14').Value = "'68.775.93"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('E15').Value = '  +3.07%  '
$ws.Range('D16').Value = "'23.56"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = "'10.57"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = "'338.66"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = "'6.91"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('D21').Value = "'1.89"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.42%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').Value = "'66.74"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = "'3.68"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').Value = "'8.19"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.77%  '
$ws.Range('D26').Value = "'0.0₃0821"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').Value = "'7.20"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = "'427.25"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('D32').Value = "'159.25"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = "'17.85"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('D37').Value = "'4.41"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.07%  '
$ws.Range('D38').Value = "'0.297"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = "'1.48"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').Value = "'2.05"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('D42').Value = "'3.37"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('D43').Value = "'130.44"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('D45').Value = "'0.483"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.49%  '
$ws.Range('D46').Value = "'0.564"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('D47').Value = "'0.0921"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.55%  '
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').Value = "'4.96"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('D51').Value = "'16.82"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.04%  '
